# [22.1.2024][Add Error invoice Scenarios]
# Adds two new test-case blocks ("testErrorClearanceInvoice" and
# "testErrorReportingInvoice") to the TestCases sheet, mirroring the layout
# of the existing "testWarningClearanceInvoice" / "testWarningReportingInvoice"
# blocks, and normalizes the stray "s=6" formatting on B23:B25 back to the
# common data-row style (s=3) used everywhere else in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) B23:B25 were carrying a one-off cell format (fill+border) that is
#    visually identical to the regular bordered data-row format used by
#    every other data cell (B3, B26, B27, ...). Re-apply that common
#    format so the stray style is no longer used.
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("B23:B25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Append the two new scenario blocks after the existing data
#    (rows 32-39), copying the header/data row formatting from the
#    existing "testWarningClearanceInvoice" / "testWarningReportingInvoice"
#    blocks (rows 22/23 and 28/29) so new cells land on the same shared
#    style records instead of minting new ones.
# ---------------------------------------------------------------------

# Row 32 - header row for testErrorClearanceInvoice
$ws.Range("A22:D22").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A32").Value = "testErrorClearanceInvoice"
$ws.Range("C32").Value = "invoiceType"
$ws.Range("D32").Value = "invoiceFileName"

# Rows 33-35 - data rows for testErrorClearanceInvoice
$ws.Range("A23:D23").Copy()
$ws.Range("A33:D35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A33").Value = "testErrorClearanceInvoice"
$ws.Range("B33").Value = $true
$ws.Range("C33").Value = "STANDARDNOTE"
$ws.Range("D33").Value = "BR-KSA-F-06-C17.xml"

$ws.Range("A34").Value = "testErrorClearanceInvoice"
$ws.Range("B34").Value = $true
$ws.Range("C34").Value = "STANDARDCREDIT"
$ws.Range("D34").Value = "BR-CL-03.xml"

$ws.Range("A35").Value = "testErrorClearanceInvoice"
$ws.Range("B35").Value = $true
$ws.Range("C35").Value = "STANDARDDEBIT"
$ws.Range("D35").Value = "BR-KSA-EN16931-08.xml"

# Row 36 - header row for testErrorReportingInvoice
$ws.Range("A28:D28").Copy()
$ws.Range("A36:D36").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A36").Value = "testErrorReportingInvoice"
$ws.Range("C36").Value = "invoiceType"
$ws.Range("D36").Value = "invoiceFileName"

# Rows 37-39 - data rows for testErrorReportingInvoice
$ws.Range("A29:D29").Copy()
$ws.Range("A37:D39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A37").Value = "testErrorReportingInvoice"
$ws.Range("B37").Value = $true
$ws.Range("C37").Value = "SIMPLIFIEDNOTE"
$ws.Range("D37").Value = "BR-KSA-49.xml"

$ws.Range("A38").Value = "testErrorReportingInvoice"
$ws.Range("B38").Value = $true
$ws.Range("C38").Value = "SIMPLIFIEDCREDIT"
$ws.Range("D38").Value = "BR-KSA-72.xml"

$ws.Range("A39").Value = "testErrorReportingInvoice"
$ws.Range("B39").Value = $true
$ws.Range("C39").Value = "SIMPLIFIEDDEBIT"
$ws.Range("D39").Value = "BR-KSA-06.xml"

# ---------------------------------------------------------------------
# 3) Move the view so the newly added rows are visible, matching the
#    updated selection/scroll position of the sheet.
# ---------------------------------------------------------------------
$ws.Range("B40").Select()
